$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "29.372.70"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "1.881.68"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'0.7134"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").Value = "'242.45"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.08046"
$ws.Range("E8").Value = "  +3.78%  "
$ws.Range("D9").Value = "'0.3131"
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("D10").Value = "'25.28"
$ws.Range("E10").Value = "  +1.40%  "
$ws.Range("D11").Value = "'0.08342"
$ws.Range("E11").Value = "  -2.31%  "
$ws.Range("D12").Value = "1.915.39"
$ws.Range("E12").Value = "  +2.08%  "
$ws.Range("D13").Value = "'5.250"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").Value = "'0.7194"
$ws.Range("E14").Value = "  +1.22%  "
$ws.Range("D15").Value = "'94.05"
$ws.Range("E15").Value = "  +2.81%  "
$ws.Range("D16").Value = "'6.337"
$ws.Range("E16").Value = "  +5.45%  "
$ws.Range("D17").Value = "'0.000008584"
$ws.Range("E17").Value = "  +4.69%  "
$ws.Range("D18").Value = "29.396.23"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").Value = "'242.61"
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("D20").Value = "2.142.93"
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").Value = "'13.26"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D23").Value = "'7.869"
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "'0.1593"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("D26").Value = "'163.59"
$ws.Range("E26").Value = "  +0.42%  "
$ws.Range("D27").Value = "'9.086"
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("D28").Value = "'18.63"
$ws.Range("E28").Value = "  +0.68%  "
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("D30").Value = "'4.421"
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("D31").Value = "'4.329"
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("D32").Value = "'1.199"
$ws.Range("E32").Value = "  -6.41%  "
$ws.Range("D33").Value = "'0.05385"
$ws.Range("E33").Value = "  +2.43%  "
$ws.Range("D34").Value = "'1.949"
$ws.Range("E34").Value = "  +0.80%  "
$ws.Range("D35").Value = "'1.183"
$ws.Range("E35").Value = "  +0.67%  "
$ws.Range("D36").Value = "'0.7508"
$ws.Range("E36").Value = "  +0.63%  "
$ws.Range("D37").Value = "'2.697"
$ws.Range("D38").Value = "'0.01892"
$ws.Range("E38").Value = "  +1.36%  "
$ws.Range("D39").Value = "1.288.55"
$ws.Range("E39").Value = "  +9.00%  "
$ws.Range("D40").Value = "'2.747"
$ws.Range("E40").Value = "  +1.13%  "
$ws.Range("D41").Value = "'6.600"
$ws.Range("E41").Value = "  +3.28%  "
$ws.Range("D42").Value = "'0.9178"
$ws.Range("E42").Value = "  +3.50%  "
$ws.Range("D43").Value = "'75.04"
$ws.Range("E43").Value = "  +2.94%  "
$ws.Range("D44").Value = "'111.90"
$ws.Range("E44").Value = "  +5.16%  "
$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("E46").Value = "  +5.41%  "
$ws.Range("D47").Value = "2.038.89"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("D48").Value = "'1.813"
$ws.Range("E48").Value = "  +0.18%  "
$ws.Range("D49").Value = "'0.5220"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").Value = "'9.534"
$ws.Range("E50").Value = "  +1.48%  "
$ws.Range("D51").Value = "'0.4395"
$ws.Range("E51").Value = "  +1.82%  "
